# Fix broken markdown link in the F3 "Speaker(s)" cell (missing "(" before the
# Mark Prell url), then restore the view: scroll back to the top of the sheet
# and move the selection to F3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "[Wen You](https://dataifa.github.io/difa-project/comingsoon.html), [Nichole Szembrot](https://dataifa.github.io/difa-project/comingsoon.html), [Mark Prell](https://dataifa.github.io/difa-project/comingsoon.html), [Bruce Weinberg](https://dataifa.github.io/difa-project/comingsoon.html)"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F3").Select()
